# Update stackcachesize.xlsx: new GETCONSTARRAY results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update raw benchmark data values (columns H, L, O across data rows) ---
$ws.Range("H5").Value = 31.5
$ws.Range("L5").Value = 216.3
$ws.Range("O5").Value = 132.5
$ws.Range("L6").Value = 20.5
$ws.Range("O6").Value = 8.8000000000000007
$ws.Range("H7").Value = 17.899999999999999
$ws.Range("L7").Value = 125.1
$ws.Range("O7").Value = 80.599999999999994
$ws.Range("H8").Value = 1.2
$ws.Range("L8").Value = 16.600000000000001
$ws.Range("O8").Value = 2.2999999999999998
$ws.Range("H9").Value = 10.4
$ws.Range("L9").Value = 54.2
$ws.Range("O9").Value = 40.700000000000003
$ws.Range("H11").Value = 29.2
$ws.Range("L11").Value = 203.2
$ws.Range("O11").Value = 115.8
$ws.Range("L12").Value = 20.5
$ws.Range("H13").Value = 15.3
$ws.Range("L13").Value = 112.8
$ws.Range("O13").Value = 64.099999999999994
$ws.Range("H14").Value = 1.5
$ws.Range("L14").Value = 16.600000000000001
$ws.Range("O14").Value = 3.7
$ws.Range("H15").Value = 10.4
$ws.Range("L15").Value = 53.3
$ws.Range("O15").Value = 40.1
$ws.Range("H17").Value = 24
$ws.Range("L17").Value = 185.8
$ws.Range("O17").Value = 101.8
$ws.Range("L18").Value = 20.5
$ws.Range("O18").Value = 6.8
$ws.Range("H19").Value = 9.4
$ws.Range("L19").Value = 101
$ws.Range("O19").Value = 50.5
$ws.Range("H20").Value = 2.2000000000000002
$ws.Range("L20").Value = 15.2
$ws.Range("O20").Value = 4.7
$ws.Range("H21").Value = 10.4
$ws.Range("L21").Value = 49
$ws.Range("O21").Value = 39.799999999999997
$ws.Range("H23").Value = 21.8
$ws.Range("L23").Value = 168.4
$ws.Range("O23").Value = 91.6
$ws.Range("L24").Value = 20.5
$ws.Range("O24").Value = 6.5
$ws.Range("H25").Value = 6.5
$ws.Range("L25").Value = 85
$ws.Range("O25").Value = 40.4
$ws.Range("H26").Value = 2.8
$ws.Range("L26").Value = 13.9
$ws.Range("O26").Value = 4.9000000000000004
$ws.Range("H27").Value = 10.4
$ws.Range("L27").Value = 49
$ws.Range("O27").Value = 39.799999999999997
$ws.Range("H29").Value = 20.9
$ws.Range("L29").Value = 161.19999999999999
$ws.Range("O29").Value = 83.4
$ws.Range("L30").Value = 20.5
$ws.Range("O30").Value = 5.9
$ws.Range("H31").Value = 5.5
$ws.Range("L31").Value = 76.900000000000006
$ws.Range("O31").Value = 31.8
$ws.Range("H32").Value = 2.9
$ws.Range("L32").Value = 14.9
$ws.Range("O32").Value = 5.9
$ws.Range("H33").Value = 10.4
$ws.Range("L33").Value = 49
$ws.Range("O33").Value = 39.799999999999997
$ws.Range("H35").Value = 19.5
$ws.Range("L35").Value = 156.30000000000001
$ws.Range("O35").Value = 77.900000000000006
$ws.Range("L36").Value = 20.399999999999999
$ws.Range("O36").Value = 5.6
$ws.Range("H37").Value = 4.2
$ws.Range("L37").Value = 72.099999999999994
$ws.Range("O37").Value = 25.9
$ws.Range("L38").Value = 14.9
$ws.Range("O38").Value = 6.7
$ws.Range("H39").Value = 10.4
$ws.Range("L39").Value = 49
$ws.Range("O39").Value = 39.799999999999997
$ws.Range("H41").Value = 18.399999999999999
$ws.Range("L41").Value = 156.4
$ws.Range("O41").Value = 72.099999999999994
$ws.Range("L42").Value = 20.399999999999999
$ws.Range("O42").Value = 5.0999999999999996
$ws.Range("H43").Value = 3.4
$ws.Range("L43").Value = 72.099999999999994
$ws.Range("O43").Value = 20.399999999999999
$ws.Range("H44").Value = 2.6
$ws.Range("L44").Value = 14.9
$ws.Range("O44").Value = 6.8
$ws.Range("H45").Value = 10.4
$ws.Range("L45").Value = 49
$ws.Range("O45").Value = 39.799999999999997
$ws.Range("L47").Value = 156.30000000000001
$ws.Range("O47").Value = 70.5
$ws.Range("L48").Value = 20.399999999999999
$ws.Range("O48").Value = 4.5
$ws.Range("H49").Value = 2.4
$ws.Range("L49").Value = 72
$ws.Range("O49").Value = 19.3
$ws.Range("L50").Value = 14.9
$ws.Range("O50").Value = 6.9
$ws.Range("H51").Value = 10.4
$ws.Range("L51").Value = 49
$ws.Range("O51").Value = 39.799999999999997

# --- Add "UPDATED 20180301" marker cell with red fill, next to row 11 ---
$ws.Range("R11").Value = "UPDATED 20180301"
$ws.Range("R11:S11").Interior.Color = 255

# --- Update the active selection to Q22 (also clears the frozen topLeftCell) ---
$null = $ws.Range("Q22").Select()
